$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B" = 1.02; "C" = 1.070955003893757; "D" = 1.07606047135903; "E" = 1.065108339823105; "F" = 1.084965789169289; "I" = 1.062660197921191; "J" = 1.07588204269338; "K" = 1.078745398953191; "L" = 1.067822468811517; "M" = 1.087627441953887 }
    3 = @{ "B" = 1.02; "C" = 1.072325022927828; "D" = 1.077198851535018; "E" = 1.066316203070301; "F" = 1.086236525146257; "I" = 1.063201114102706; "J" = 1.076907821865162; "K" = 1.07970042639513; "L" = 1.068844685466319; "M" = 1.088716194643674 }
    4 = @{ "B" = 1.02; "C" = 1.073210065510805; "D" = 1.077934222002555; "E" = 1.067095971961824; "F" = 1.087057808314186; "I" = 1.063549131949797; "J" = 1.077569620331844; "K" = 1.08031656859472; "L" = 1.069503760708312; "M" = 1.089419159242752 }
    5 = @{ "B" = 1.02; "C" = 1.073581795873476; "D" = 1.078243079671593; "E" = 1.067423361206657; "F" = 1.087402848627433; "I" = 1.063694964593371; "J" = 1.077847378251784; "K" = 1.080575162022376; "L" = 1.069780274036463; "M" = 1.089714323291804 }
    6 = @{ "B" = 1.02; "C" = 1.073644191180412; "D" = 1.078294921227861; "E" = 1.067478306494346; "F" = 1.087460769228625; "I" = 1.063719422812495; "J" = 1.077893988058016; "K" = 1.080618555736188; "L" = 1.069826669030598; "M" = 1.089763861515427 }
    7 = @{ "B" = 1.02; "C" = 1.073215033924773; "D" = 1.077938350116273; "E" = 1.067100348218803; "F" = 1.087062419646488; "I" = 1.063551082432213; "J" = 1.07757333355714; "K" = 1.080320025628202; "L" = 1.069507457694837; "M" = 1.089423104657352 }
    8 = @{ "B" = 1.02; "C" = 1.071418312319542; "D" = 1.076445451072666; "E" = 1.065516918566541; "F" = 1.085395443549038; "I" = 1.062843416922727; "J" = 1.076229115886993; "K" = 1.079068535129871; "L" = 1.068168424792614; "M" = 1.08799571066357 }
    9 = @{ "B" = 1.02; "C" = 1.068240877960936; "D" = 1.07380510375021; "E" = 1.062712699032933; "F" = 1.082450397870893; "I" = 1.061581052946451; "J" = 1.073845301376757; "K" = 1.076849095439671; "L" = 1.065790530059978; "M" = 1.085468539376886 }
    10 = @{ "B" = 1.02; "C" = 1.066114562561245; "D" = 1.072038104032135; "E" = 1.060833469181776; "F" = 1.080481618091812; "I" = 1.060728979362916; "J" = 1.072245633443682; "K" = 1.075359698263003; "L" = 1.064192625556582; "M" = 1.083775470831689 }
    11 = @{ "B" = 1.02; "C" = 1.065191857140124; "D" = 1.071271308054801; "E" = 1.060017355716843; "F" = 1.0796277643609; "I" = 1.060357494890998; "J" = 1.071550415009402; "K" = 1.074712398066149; "L" = 1.063497648373693; "M" = 1.083040326337599 }
    12 = @{ "B" = 1.02; "C" = 1.064848815960367; "D" = 1.070986229187464; "E" = 1.059713848654372; "F" = 1.079310395122609; "I" = 1.060219125481587; "J" = 1.071291790966704; "K" = 1.074471599150723; "L" = 1.063239035076239; "M" = 1.082766949878216 }
    13 = @{ "B" = 1.02; "C" = 1.064922413481802; "D" = 1.071047391273714; "E" = 1.059778968623903; "F" = 1.079378481494882; "I" = 1.060248823598408; "J" = 1.071347284393501; "K" = 1.074523267834076; "L" = 1.063294529773325; "M" = 1.082825604219767 }
    14 = @{ "B" = 1.02; "C" = 1.065163507560497; "D" = 1.071247748614001; "E" = 1.059992275235784; "F" = 1.079601534836914; "I" = 1.06034606507555; "J" = 1.071529045032601; "K" = 1.074692500982908; "L" = 1.06347628089622; "M" = 1.083017735331014 }
    15 = @{ "B" = 1.02; "C" = 1.065312012639562; "D" = 1.071371161210006; "E" = 1.060123651763262; "F" = 1.079738937369988; "I" = 1.060405927838091; "J" = 1.071640982051737; "K" = 1.074796722885693; "L" = 1.063588201589801; "M" = 1.083136072269438 }
    16 = @{ "B" = 1.02; "C" = 1.06617575678823; "D" = 1.072088958048185; "E" = 1.060887580946484; "F" = 1.080538256426921; "I" = 1.060753579959093; "J" = 1.072291718563556; "K" = 1.075402606844057; "L" = 1.06423868362709; "M" = 1.083824216624742 }
    17 = @{ "B" = 1.02; "C" = 1.066717020926362; "D" = 1.072538761295078; "E" = 1.061366127927159; "F" = 1.081039280234531; "I" = 1.060970973021195; "J" = 1.072699221288857; "K" = 1.075782020742053; "L" = 1.064645886608839; "M" = 1.0842553228807 }
    18 = @{ "B" = 1.02; "C" = 1.067032538972875; "D" = 1.072800962816099; "E" = 1.061645025683792; "F" = 1.081331388203751; "I" = 1.061097530736168; "J" = 1.072936664916769; "K" = 1.076003097021784; "L" = 1.064883105031118; "M" = 1.084506583918807 }
    19 = @{ "B" = 1.02; "C" = 1.067140090097916; "D" = 1.07289033959199; "E" = 1.06174008369116; "F" = 1.081430967469712; "I" = 1.061140642361984; "J" = 1.07301758553155; "K" = 1.076078439556293; "L" = 1.064963940324808; "M" = 1.084592224365453 }
    20 = @{ "B" = 1.02; "C" = 1.066658968338158; "D" = 1.072490518366787; "E" = 1.061314808265241; "F" = 1.080985538687892; "I" = 1.060947674063613; "J" = 1.072655525570932; "K" = 1.075741336989925; "L" = 1.064602228263493; "M" = 1.084209089567401 }
    21 = @{ "B" = 1.02; "C" = 1.065092519905397; "D" = 1.071188755499872; "E" = 1.059929471915464; "F" = 1.07953585702909; "I" = 1.060317440498178; "J" = 1.071475531837091; "K" = 1.074642676094226; "L" = 1.063422772683249; "M" = 1.082961166137733 }
    22 = @{ "B" = 1.02; "C" = 1.064105849991593; "D" = 1.070368797687291; "E" = 1.059056334123167; "F" = 1.078623167075516; "I" = 1.059918966976646; "J" = 1.070731369727903; "K" = 1.073949802812788; "L" = 1.06267849244174; "M" = 1.082174744746366 }
    23 = @{ "B" = 1.02; "C" = 1.06462907370791; "D" = 1.07080361564111; "E" = 1.059519404271737; "F" = 1.079107118631064; "I" = 1.060130416925207; "J" = 1.071126079601026; "K" = 1.07431730893981; "L" = 1.063073308340393; "M" = 1.082591814258137 }
    24 = @{ "B" = 1.02; "C" = 1.06668520039716; "D" = 1.072512317766928; "E" = 1.061337998123501; "F" = 1.081009822583139; "I" = 1.060958202613276; "J" = 1.072675270544139; "K" = 1.075759720934846; "L" = 1.064621956505238; "M" = 1.084229981016673 }
    25 = @{ "B" = 1.02; "C" = 1.069063707975155; "D" = 1.074488870017843; "E" = 1.063439351369603; "F" = 1.083212695830281; "I" = 1.061909242346292; "J" = 1.074463396906063; "K" = 1.077424577174966; "L" = 1.066407478462817; "M" = 1.086123314369631 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}

Write-Host "Applied 380 kV case updates"